# Fruta / hortaliza, semanal
# New weekly price record inserted as row 148 ("Crimpson Seedless"),
# pushing the former rows 148-155 down to rows 149-156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 148; this shifts rows 148:155
# down to 149:156 (including their formatting) and grows the sheet
# dimension to A1:T156 automatically.
$ws.Rows("148:148").Insert()

# Populate the newly inserted row 148 with this week's record.
$ws.Range("A148").Value = 11
$ws.Range("B148").Value = "Vega Monumental Concepción"
$ws.Range("C148").Value = "Bíobío"
$ws.Range("D148").Value = 44714
$ws.Range("E148").Value = 8
$ws.Range("F148").Value = "Fruta"
$ws.Range("G148").Value = 100109
$ws.Range("H148").Value = "Uva"
$ws.Range("I148").Value = 100109001
$ws.Range("J148").Value = "Uva"
$ws.Range("K148").Value = "Crimpson Seedless"
$ws.Range("L148").Value = "Primera"
$ws.Range("M148").Value = 100
$ws.Range("N148").Value = 9000
$ws.Range("O148").Value = 10000
$ws.Range("P148").Value = 9500
$ws.Range("Q148").Value = "`$/bandeja 18 kilos"
$ws.Range("R148").Value = "Región de O'Higgins"
$ws.Range("S148").Value = 528
$ws.Range("T148").Value = 18
